$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Titles / headers (new rows 1-4, column D) ---------------------------
$ws.Range("D1").Value = "NR Finance Mexico"
$ws.Range("D2").Value = "PROCOTIZA"
$ws.Range("D3").Value = "Certificacion de usuarios 2024"
$ws.Range("D4").Value = "Reporte de usuarios"

$titles = $ws.Range("D1:D4")
$titles.Font.Name = "Calibri"
$titles.Font.Bold = $true
$titles.Font.Size = 16
$titles.HorizontalAlignment = -4108

# --- Borders around the existing data table (rows 5-8, cols A-H) ---------
$data = $ws.Range("A5:H8")
$data.Borders.Color = 0
$data.Borders.LineStyle = 1

# --- Column widths (manual "auto-fit") ------------------------------------
$ws.Columns("A").ColumnWidth = 11.282054
$ws.Columns("B").ColumnWidth = 37.139196
$ws.Columns("C").ColumnWidth = 8.424911
$ws.Columns("D").ColumnWidth = 40.424911
$ws.Columns("E").ColumnWidth = 21.567768
$ws.Columns("F:G").ColumnWidth = 24.139196
$ws.Columns("H").ColumnWidth = 21.853482
